$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 146, shifting existing rows 146-177 down to 147-178
$ws.Rows.Item(146).Insert()

# Populate the new row 146 with the new record's data
$ws.Cells.Item(146, 1).Value = 4
$ws.Cells.Item(146, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(146, 3).Value = "Los Lagos"
$ws.Cells.Item(146, 4).Value = 44637
$ws.Cells.Item(146, 5).Value = 10
$ws.Cells.Item(146, 6).Value = "Fruta"
$ws.Cells.Item(146, 7).Value = 100108
$ws.Cells.Item(146, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(146, 9).Value = 100108002
$ws.Cells.Item(146, 10).Value = "Mango"
$ws.Cells.Item(146, 11).Value = "Sin especificar"
$ws.Cells.Item(146, 12).Value = "Primera"
$ws.Cells.Item(146, 13).Value = 80
$ws.Cells.Item(146, 14).Value = 7000
$ws.Cells.Item(146, 15).Value = 7500
$ws.Cells.Item(146, 16).Value = 7250
$ws.Cells.Item(146, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(146, 18).Value = "Perú"
$ws.Cells.Item(146, 19).Value = 1812
$ws.Cells.Item(146, 20).Value = 4
